# commented out HW UI update
# The hardware UI code that normally pushes live Ticket Sales (Q) and
# Embarking (R) counts into the sheet was commented out, so these cells
# now reflect the last simulated/random values instead of the live
# hardware-driven ones. Apply those updated values here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> [TicketSales(Q), Embarking(R)] new values.
# $null means "leave this cell unchanged".
$updates = @{
    3   = @(4, 3)
    10  = @(76, 66)
    17  = @(40, 20)
    23  = @(30, 5)
    32  = @(42, 11)
    40  = @(13, 11)
    49  = @(32, 28)
    58  = @(63, 27)
    66  = @(56, 52)
    74  = @(93, 79)
    78  = @(82, 56)
    89  = @(69, 50)
    97  = @(10, 9)
    106 = @($null, 27)
    115 = @(1, 1)
    124 = @(81, 11)
    133 = @(17, 7)
    142 = @(70, 53)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $q = $vals[0]
    $r = $vals[1]

    if ($null -ne $q) {
        $ws.Cells.Item($row, 17).Value = $q   # Column Q = Ticket Sales
    }
    if ($null -ne $r) {
        $ws.Cells.Item($row, 18).Value = $r   # Column R = Embarking
    }
}
